# Apply the cryptos list refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "20.80", "157.00").
# Force the whole Price column to Text format first so Excel keeps them as
# literal strings (preserving trailing zeros / multi-dot "thousand" prices)
# instead of silently coercing them into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.121.07"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.882.97"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "313.62"

$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "0.5081"
$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("D8").Value = "0.3868"
$ws.Range("E8").Value = "  -1.34%  "

$ws.Range("D9").Value = "0.09036"
$ws.Range("E9").Value = "  -3.57%  "

$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("D11").Value = "41.72"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Value = "6.379"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "20.82"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").Value = "1.882.54"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "7.265"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("D19").Value = "0.06633"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("E20").Value = "  +2.12%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "6.128"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("D23").Value = "28.152.10"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("D25").Value = "2.268"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").Value = "2.096.36"
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("E27").Value = "  -3.05%  "

$ws.Range("D28").Value = "20.80"
$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("D29").Value = "157.00"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Value = "127.17"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("E32").Value = "  -1.81%  "

$ws.Range("D33").Value = "5.630"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").Value = "3.605"
$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("D35").Value = "9.597"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").Value = "0.06639"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "0.02413"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").Value = "1.215"
$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("D41").Value = "0.6423"
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").Value = "11.50"
$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").Value = "4.923"
$ws.Range("E43").Value = "  -1.53%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6058"
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.21"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "1.277"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "3.672"
$ws.Range("E48").Value = "  -1.18%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.249"
$ws.Range("E49").Value = "  +6.04%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "2.006"
$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "121.51"
$ws.Range("E51").Value = "  -0.60%  "
